$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2
$ws.Range("F5").Value = 2
$ws.Range("F8").Value = 0
$ws.Range("F12").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("F15").Value = -1
$ws.Range("F17").Value = 2
$ws.Range("F20").Value = -1
